# Germanize the "2012 Sales" workbook: translate headers, region names and
# month names to German, swap the Order Amount column from USD to EUR
# accounting format, rename the table columns, and tweak a couple of
# cosmetic details (column width, selected cell).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the worksheet tab.
$ws.Name = "2012 Verkäufe"

# 2) Translate the header row (B1 "Region" is already identical in German).
$ws.Range("A1").Value = "Verkäufer"
$ws.Range("C1").Value = "Konto"
$ws.Range("E1").Value = "Monat"

# 3) Translate the Region column (B2:B40) values, whole-cell matches only.
$regionCol = $ws.Range("B2:B40")
$regionCol.Replace("East", "Osten", 1)
$regionCol.Replace("West", "Westen", 1)
$regionCol.Replace("North", "Norden", 1)
$regionCol.Replace("South", "Süden", 1)

# 4) Translate the Month column (E2:E40) values, whole-cell matches only.
$monthCol = $ws.Range("E2:E40")
$monthCol.Replace("January", "Januar", 1)
$monthCol.Replace("February", "Februar", 1)
$monthCol.Replace("March", "März", 1)

# Order Amount header translates last.
$ws.Range("D1").Value = "Umsatz"

# 5) Switch the Order Amount column from a USD to a EUR accounting format.
$eurFormat = "_-* #,##0.00\ [$€-407]_-;\-* #,##0.00\ [$€-407]_-;_-* ""-""??\ [$€-407]_-;_-@_-"
$ws.Range("D2:D39").NumberFormat = $eurFormat
$ws.Range("D40").NumberFormat = $eurFormat

# 6) Rename the table ("Salesperson"/"Account"/"Order Amount"/"Month" ->
#    German); the header cells already carry the new text, this keeps the
#    ListObject's column metadata (used for dxf-based column formatting) in
#    sync.
$tbl = $ws.ListObjects.Item(1)
$tbl.ListColumns.Item(1).Name = "Verkäufer"
$tbl.ListColumns.Item(3).Name = "Konto"
$tbl.ListColumns.Item(4).Name = "Umsatz"
$tbl.ListColumns.Item(5).Name = "Monat"

# 7) Cosmetic tweaks: narrower Order Amount column (shorter German header),
#    and leave the cursor on E5 like the edited workbook.
$ws.Columns.Item(4).ColumnWidth = 12.5
$ws.Range("E5").Select()
